$d = $word.ActiveDocument

function ReplaceText($old, $new) {
    $null = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

# ---------------------------------------------------------------------------
# 1. Topic paragraph: expand "Cobra Bank App…" into the full description and
#    wrap "that allows users to" in the _Hlk73462928 bookmark.
# ---------------------------------------------------------------------------
ReplaceText "Cobra Bank App…" "Cobra Bank App is a simple web application that allows users to view their accounts, transfer, deposit and withdraw money from their account. "

$bmRange = $d.Content
$null = $bmRange.Find.Execute("that allows users to", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$d.Bookmarks.Add("_Hlk73462928", $bmRange)

# ---------------------------------------------------------------------------
# 2. Requirement #2: append " and log out"
# ---------------------------------------------------------------------------
ReplaceText "This system shall allow users to log in" "This system shall allow users to log in and log out"

# ---------------------------------------------------------------------------
# 3. Cascade the existing requirement descriptions down a slot, freeing up
#    requirement #3 for the new "register" requirement. Must run top to
#    bottom so each search target still exists (unmodified) when searched.
# ---------------------------------------------------------------------------
ReplaceText "This system shall log user access attempts" "The system shall allow new users to register"
ReplaceText "This system shall allow users to make deposits" "This system shall log user access attempts"
ReplaceText "This system shall allow users to make withdrawals" "This system shall allow users to make deposits"
ReplaceText "This system shall allow users to view their balance" "This system shall allow users to make withdrawals"
ReplaceText "This system shall allow users to transfer to other users" "This system shall allow users to view their balance"

# ---------------------------------------------------------------------------
# 4. Requirement table: fill in the previously-empty row 8 description and
#    append five brand-new requirement rows (9-13).
# ---------------------------------------------------------------------------
$reqTable = $d.Tables.Item(1)
$reqTable.Rows.Item($reqTable.Rows.Count).Cells.Item(2).Range.Text = "This system shall only allow a savings and checking account"

$newRequirements = @(
    "This system shall allow users to transfer balances between account types",
    "This system shall allow accounts to collect interest",
    "This system shall prompt users for verification before executing actions",
    "This system shall alert users for invalid transactions to include overdrawing accounts",
    "This system shall not prevent account overdraft"
)

# Existing table has a header row plus rows labelled 1-8, so the next label is 9.
$reqNumber = $reqTable.Rows.Count - 1
foreach ($reqText in $newRequirements) {
    $reqNumber = $reqNumber + 1
    $row = $reqTable.Rows.Add()
    $row.Cells.Item(1).Range.Text = [string]$reqNumber
    $row.Cells.Item(2).Range.Text = $reqText
}

# ---------------------------------------------------------------------------
# 5. Revision history table: log the 6/1 retry session.
# ---------------------------------------------------------------------------
$revTable = $d.Tables.Item(2)

$row = $revTable.Rows.Add()
$row.Cells.Item(1).Range.Text = "6/1"
$row.Cells.Item(2).Range.Text = "Keith"
$row.Cells.Item(3).Range.Text = "Added 2 additional items, phrasing to the topic"

$row = $revTable.Rows.Add()
$row.Cells.Item(1).Range.Text = "6/1"
$row.Cells.Item(2).Range.Text = "All"
$row.Cells.Item(3).Range.Text = "Additional requirements and revision over zoom call "
